# Apply Seraph_Profits profit-recalculation updates across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4199
$ws.Range("I32").Value = 3998.6667
$ws.Range("J32").Value = 4800
$ws.Range("K32").Value = 3998.6667
$ws.Range("L32").Value = 4800
$ws.Range("M32").Value = -3672.6667
$ws.Range("N32").Value = -5452
$ws.Range("H33").Value = 211.76923
$ws.Range("J33").Value = 445.6
$ws.Range("L33").Value = 445.6
$ws.Range("N33").Value = -903.6
$ws.Range("H51").Value = 9000
$ws.Range("I51").Value = 9000
$ws.Range("K51").Value = 9000
$ws.Range("M51").Value = -8516
$ws.Range("H112").Value = 2391.647
$ws.Range("J112").Value = 2739.9285
$ws.Range("L112").Value = 8219.7855
$ws.Range("N112").Value = -10435.7855
$ws.Range("H125").Value = 2500
$ws.Range("J125").Value = 4000
$ws.Range("L125").Value = 36000
$ws.Range("N125").Value = -40920
$ws.Range("H132").Value = 1383.4857
$ws.Range("I132").Value = 1326.4242
$ws.Range("K132").Value = 3979.2726
$ws.Range("M132").Value = -1449.2726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 616.6
$ws.Range("I2").Value = 419.66666
$ws.Range("J2").Value = 1404.3334
$ws.Range("K2").Value = 419.66666
$ws.Range("L2").Value = 1404.3334
$ws.Range("M2").Value = -306.66666
$ws.Range("N2").Value = -1630.3334
$ws.Range("H32").Value = 13862.873
$ws.Range("I32").Value = 5574.4736
$ws.Range("K32").Value = 5574.4736
$ws.Range("M32").Value = -5287.4736
$ws.Range("H45").Value = 1656
$ws.Range("I45").Value = 1656
$ws.Range("K45").Value = 1656
$ws.Range("M45").Value = -1279
$ws.Range("H110").Value = 5410.2856
$ws.Range("I110").Value = 6941.8
$ws.Range("J110").Value = 1581.5
$ws.Range("K110").Value = 6941.8
$ws.Range("L110").Value = 1581.5
$ws.Range("M110").Value = -4896.8
$ws.Range("N110").Value = -5671.5
$ws.Range("H116").Value = 616.6
$ws.Range("I116").Value = 419.66666
$ws.Range("J116").Value = 1404.3334
$ws.Range("K116").Value = 419.66666
$ws.Range("L116").Value = 1404.3334
$ws.Range("M116").Value = 1874.33334
$ws.Range("N116").Value = -5992.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 616.6
$ws.Range("I3").Value = 419.66666
$ws.Range("J3").Value = 1404.3334
$ws.Range("K3").Value = 419.66666
$ws.Range("L3").Value = 1404.3334
$ws.Range("M3").Value = -305.66666
$ws.Range("N3").Value = -1632.3334
$ws.Range("H80").Value = 585.8182
$ws.Range("I80").Value = 522.2857
$ws.Range("J80").Value = 697
$ws.Range("K80").Value = 522.2857
$ws.Range("L80").Value = 697
$ws.Range("M80").Value = 475.7143
$ws.Range("N80").Value = -2693
$ws.Range("H83").Value = 585.8182
$ws.Range("I83").Value = 522.2857
$ws.Range("J83").Value = 697
$ws.Range("K83").Value = 2611.4285
$ws.Range("L83").Value = 3485
$ws.Range("M83").Value = 2380.5715
$ws.Range("N83").Value = -13469
$ws.Range("H105").Value = 3582.1785
$ws.Range("J105").Value = 6041
$ws.Range("L105").Value = 6041
$ws.Range("N105").Value = -9535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 940
$ws.Range("H52").Value = 99749.5
$ws.Range("J52").Value = 99749.5
$ws.Range("L52").Value = 99749.5
$ws.Range("N52").Value = -100337.5
$ws.Range("H94").Value = 1750
$ws.Range("J94").Value = 1750
$ws.Range("L94").Value = 1750
$ws.Range("N94").Value = -2652
$ws.Range("H134").Value = 2291.1482
$ws.Range("I134").Value = 1422.45
$ws.Range("K134").Value = 4267.35
$ws.Range("M134").Value = -1732.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1235.1428
$ws.Range("I68").Value = 1199.6666
$ws.Range("J68").Value = 1261.75
$ws.Range("K68").Value = 3598.9998
$ws.Range("L68").Value = 3785.25
$ws.Range("M68").Value = -2787.9998
$ws.Range("N68").Value = -5407.25
$ws.Range("H71").Value = 1235.1428
$ws.Range("I71").Value = 1199.6666
$ws.Range("J71").Value = 1261.75
$ws.Range("K71").Value = 10796.9994
$ws.Range("L71").Value = 11355.75
$ws.Range("M71").Value = -6740.999400000001
$ws.Range("N71").Value = -19467.75
$ws.Range("H80").Value = 4596
$ws.Range("I80").Value = 1794
$ws.Range("J80").Value = 5997
$ws.Range("K80").Value = 5382
$ws.Range("L80").Value = 17991
$ws.Range("M80").Value = -4446
$ws.Range("N80").Value = -19863
$ws.Range("H83").Value = 4596
$ws.Range("I83").Value = 1794
$ws.Range("J83").Value = 5997
$ws.Range("K83").Value = 16146
$ws.Range("L83").Value = 53973
$ws.Range("M83").Value = -11466
$ws.Range("N83").Value = -63333
$ws.Range("H100").Value = 6000
$ws.Range("J100").Value = 6000
$ws.Range("L100").Value = 18000
$ws.Range("N100").Value = -19622
$ws.Range("H103").Value = 1000
$ws.Range("I103").Value = 1000
$ws.Range("K103").Value = 3000
$ws.Range("M103").Value = -2121
$ws.Range("H121").Value = 12309
$ws.Range("I121").Value = 574
$ws.Range("K121").Value = 1722
$ws.Range("M121").Value = -412
$ws.Range("H122").Value = 824.75
$ws.Range("I122").Value = 599.6667
$ws.Range("K122").Value = 5397.0003
$ws.Range("M122").Value = -2947.0003
$ws.Range("H132").Value = 1536.5
$ws.Range("I132").Value = 1560
$ws.Range("J132").Value = 1513
$ws.Range("K132").Value = 14040
$ws.Range("L132").Value = 13617
$ws.Range("M132").Value = -11510
$ws.Range("N132").Value = -18677

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H96").Value = 50173.668
$ws.Range("J96").Value = 50173.668
$ws.Range("L96").Value = 50173.668
$ws.Range("N96").Value = -55665.668
$ws.Range("H98").Value = 26048.25
$ws.Range("J98").Value = 26048.25
$ws.Range("L98").Value = 26048.25
$ws.Range("N98").Value = -32038.25
$ws.Range("H99").Value = 20194
$ws.Range("I99").Value = 11235.5
$ws.Range("J99").Value = 26166.334
$ws.Range("K99").Value = 11235.5
$ws.Range("L99").Value = 26166.334
$ws.Range("M99").Value = -8989.5
$ws.Range("N99").Value = -30658.334
$ws.Range("H100").Value = 44996.5
$ws.Range("J100").Value = 44996.5
$ws.Range("L100").Value = 44996.5
$ws.Range("N100").Value = -47160.5
$ws.Range("H101").Value = 21314.666
$ws.Range("J101").Value = 21314.666
$ws.Range("L101").Value = 21314.666
$ws.Range("N101").Value = -27804.666
$ws.Range("H102").Value = 1569.742
$ws.Range("I102").Value = 348.68182
$ws.Range("J102").Value = 4554.5557
$ws.Range("K102").Value = 348.68182
$ws.Range("L102").Value = 4554.5557
$ws.Range("M102").Value = 1273.31818
$ws.Range("N102").Value = -7798.5557
$ws.Range("H105").Value = 42962.332
$ws.Range("J105").Value = 42962.332
$ws.Range("L105").Value = 42962.332
$ws.Range("N105").Value = -49950.332
$ws.Range("H107").Value = 714.7917
$ws.Range("I107").Value = 492.63635
$ws.Range("K107").Value = 492.63635
$ws.Range("M107").Value = 1427.36365
$ws.Range("H122").Value = 54019.9
$ws.Range("I122").Value = 3752.1538
$ws.Range("K122").Value = 11256.4614
$ws.Range("M122").Value = -8806.4614
$ws.Range("H126").Value = 4078.818
$ws.Range("I126").Value = 3013.6
$ws.Range("K126").Value = 9040.8
$ws.Range("M126").Value = -6570.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3835.3809
$ws.Range("I132").Value = 3447.1924
$ws.Range("J132").Value = 4466.1875
$ws.Range("K132").Value = 10341.5772
$ws.Range("L132").Value = 13398.5625
$ws.Range("M132").Value = -7811.5772
$ws.Range("N132").Value = -18458.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1478
$ws.Range("I132").Value = 1353.579
$ws.Range("K132").Value = 4060.737
$ws.Range("M132").Value = -1530.737

Write-Output "applied all updates"